$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.59"
$ws.Range("E2").Value = "'-0.12%"
$ws.Range("D3").Value = "'41.27"
$ws.Range("E3").Value = "'0.30%"
$ws.Range("D4").Value = "'5.705"
$ws.Range("E4").Value = "'-0.78%"
$ws.Range("D5").Value = "'0.08414"
$ws.Range("E5").Value = "'3.02%"
$ws.Range("D6").Value = "'8.802"
$ws.Range("E6").Value = "'0.61%"
$ws.Range("D7").Value = "'4.517"
$ws.Range("E7").Value = "'0.14%"
$ws.Range("D8").Value = "'1.982"
$ws.Range("E8").Value = "'-3.60%"
$ws.Range("D9").Value = "'2.914"
$ws.Range("E9").Value = "'-0.74%"
$ws.Range("D10").Value = "'0.9249"
$ws.Range("E10").Value = "'0.44%"
$ws.Range("D11").Value = "'0.1243"
$ws.Range("E11").Value = "'0.17%"
$ws.Range("D12").Value = "'0.1959"
$ws.Range("E12").Value = "'0.47%"
$ws.Range("D13").Value = "'0.09396"
$ws.Range("E13").Value = "'-1.09%"
$ws.Range("D14").Value = "'0.03975"
$ws.Range("E14").Value = "'8.76%"
$ws.Range("D15").Value = "'0.1064"
$ws.Range("E15").Value = "'0.84%"
$ws.Range("D16").Value = "'0.001309"
$ws.Range("E16").Value = "'1.09%"
$ws.Range("D17").Value = "'0.006114"
$ws.Range("E17").Value = "'-0.91%"
$ws.Range("D18").Value = "'3.435"
$ws.Range("E18").Value = "'1.47%"
$ws.Range("D20").Value = "'9.107"
$ws.Range("E20").Value = "'9.99%"
$ws.Range("D21").Value = "'0.1376"
$ws.Range("E21").Value = "'-2.85%"
$ws.Range("D23").Value = "'0.04428"
$ws.Range("E23").Value = "'0.00%"
$ws.Range("D24").Value = "'0.001246"
$ws.Range("E24").Value = "'-1.11%"
$ws.Range("D25").Value = "'0.004370"
$ws.Range("E25").Value = "'0.67%"
$ws.Range("D26").Value = "'0.0001197"
$ws.Range("E26").Value = "'-3.54%"
$ws.Range("D27").Value = "'0.0004007"
$ws.Range("E27").Value = "'0.36%"
$ws.Range("D39").Value = "'0.02798"
$ws.Range("E39").Value = "'0.78%"
$ws.Range("D40").Value = "'0.05530"
$ws.Range("E40").Value = "'0.28%"
$ws.Range("D41").Value = "'0.007904"
$ws.Range("E41").Value = "'3.72%"
$ws.Range("D42").Value = "'0.1435"
$ws.Range("E42").Value = "'0.74%"
$ws.Range("D43").Value = "'0.009000"
$ws.Range("E43").Value = "'-9.45%"
$ws.Range("D44").Value = "'0.002103"
$ws.Range("E44").Value = "'-1.24%"
$ws.Range("D45").Value = "'0.01013"
$ws.Range("D46").Value = "'0.00007187"
$ws.Range("E46").Value = "'6.90%"
$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.49%"
$ws.Range("D48").Value = "'0.003474"
$ws.Range("E48").Value = "'15.54%"
$ws.Range("D49").Value = "'0.002288"
$ws.Range("E49").Value = "'0.43%"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("E50").Value = "'0.49%"
$ws.Range("D51").Value = "'0.0002009"
$ws.Range("E51").Value = "'0.49%"
